# Generate Report for Handback
# Row 7 ("d90dc141-16c2-4bdb-a019-d52b402cf9ed") on both the zh-cn and de-de
# sheets gets a freshly generated handback report: the target/handback file
# info is filled in, a hyperlink is added on the "Latest Target File" cell,
# and an "Error Detail" is recorded because the handback file used was not
# built from the very latest source revision. The "Error Detail" column is
# also widened so the message is readable.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37b3aa88f02ee2d681bc23ac1236ebd135b78d0b/e2e/d90dc141-16c2-4bdb-a019-d52b402cf9ed.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/26f0dde45466f9f6f5abbc5fbf21b3588fbff5d3/e2e/d90dc141-16c2-4bdb-a019-d52b402cf9ed.md."
$currentHandbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/37b3aa88f02ee2d681bc23ac1236ebd135b78d0b/e2e/d90dc141-16c2-4bdb-a019-d52b402cf9ed.md"
$currentHandbackDisplay = "d90dc141-16c2-4bdb-a019-d52b402cf9ed.md"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Columns.Item(16).ColumnWidth = 39.17

$wsZh.Hyperlinks.Add($wsZh.Range("I7"), $currentHandbackUrl, "", "", $currentHandbackDisplay)
$wsZh.Range("J7").Value = "d90dc141-16c2-4bdb-a019-d52b402cf9ed.5c51116b37c704589240c41c2ccb94c1b6d25b2b.zh-cn.xlf"
$wsZh.Range("K7").Value = "2016-08-16 08:43:44"
$wsZh.Range("P7").Value = $errorDetail

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Columns.Item(16).ColumnWidth = 39.17

$wsDe.Hyperlinks.Add($wsDe.Range("I7"), $currentHandbackUrl, "", "", $currentHandbackDisplay)
$wsDe.Range("J7").Value = "d90dc141-16c2-4bdb-a019-d52b402cf9ed.5c51116b37c704589240c41c2ccb94c1b6d25b2b.de-de.xlf"
$wsDe.Range("K7").Value = "2016-08-16 08:43:51"
$wsDe.Range("P7").Value = $errorDetail
